$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New weekly price observation: insert a row above row 38. This pushes
# the existing rows 38-46 down to 39-47 (and carries the row's
# formatting, including the date style on column D, down with them,
# same as a normal Excel row insert).
$ws.Rows.Item(38).Insert()

# Fill in the new row 38 with the new entry. Columns A, B, C, E, F, G,
# H, I, J, K (market/product identity) are identical for this whole
# block of rows, so they're simply repeated here.
$ws.Cells.Item(38, 1).Value = 5
$ws.Cells.Item(38, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(38, 3).Value = 'Maule'
$ws.Cells.Item(38, 4).Value = 44706
$ws.Cells.Item(38, 5).Value = 7
$ws.Cells.Item(38, 6).Value = 'Fruta'
$ws.Cells.Item(38, 7).Value = 100107
$ws.Cells.Item(38, 8).Value = 'Otros'
$ws.Cells.Item(38, 9).Value = 100107001
$ws.Cells.Item(38, 10).Value = 'Caqui'
$ws.Cells.Item(38, 11).Value = 'Mankaki'
$ws.Cells.Item(38, 12).Value = 'Primera'
$ws.Cells.Item(38, 13).Value = 150
$ws.Cells.Item(38, 14).Value = 12000
$ws.Cells.Item(38, 15).Value = 12000
$ws.Cells.Item(38, 16).Value = 12000
$ws.Cells.Item(38, 17).Value = '$/caja 12 kilos empedrada'
$ws.Cells.Item(38, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(38, 19).Value = 1000
$ws.Cells.Item(38, 20).Value = 12

# Make sure the date cell keeps the same date number format style as
# the rest of the D column.
$ws.Cells.Item(38, 4).NumberFormat = $ws.Cells.Item(39, 4).NumberFormat
